$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 (id 4101): change from "철괴 / ingots" to "철덩어리 / iron"
# Keep Tooltip text the same ("다양한 걸 만들 수 있을 것 같다.")
$ws.Range("B13").Value = "철덩어리"
$ws.Range("C13").Value = "다양한 걸 만들 수 있을 것 같다."
$ws.Range("E13").Value = 20

# Row 14 (new item, id 4102): wire
$ws.Range("A14").Value = 4102
$ws.Range("B14").Value = "전선"
$ws.Range("C14").Value = "구리로 만들어진 전선."
$ws.Range("D14").Value = "wire"
$ws.Range("E14").Value = 20
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0

# Finish row 13 icon column
$ws.Range("D13").Value = "iron"

# Row 15 (new item, id 4103): gold powder
$ws.Range("A15").Value = 4103
$ws.Range("B15").Value = "금가루"
$ws.Range("C15").Value = "가치가 높은 금속의 가루"
$ws.Range("D15").Value = "gold"
$ws.Range("E15").Value = 20
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0

# Selection cosmetic change, matches the final workbook state
$ws.Range("D16").Select()
